# Applies the "Updated cryptos list" data refresh described in the commit diff.
# Numeric-looking values in column D are written with a leading apostrophe so Excel
# keeps them as literal text (matching the original inlineStr cell content) instead of
# silently coercing them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.209.99'
$ws.Range('E2').Value = '  -0.62%  '

# Row 3
$ws.Range('D3').Value = '1.826.59'
$ws.Range('E3').Value = '  -0.82%  '

# Row 4
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
$ws.Range('D5').Value = '''233.96'
$ws.Range('E5').Value = '  -2.35%  '

# Row 6
$ws.Range('D6').Value = '''0.5998'
$ws.Range('E6').Value = '  -4.26%  '

# Row 7
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
$ws.Range('D8').Value = '''0.07057'
$ws.Range('E8').Value = '  -5.20%  '

# Row 9
$ws.Range('D9').Value = '''0.2791'
$ws.Range('E9').Value = '  -3.75%  '

# Row 10
$ws.Range('D10').Value = '''23.43'
$ws.Range('E10').Value = '  -5.42%  '

# Row 11
$ws.Range('D11').Value = '''0.07641'
$ws.Range('E11').Value = '  -1.05%  '

# Row 12
$ws.Range('D12').Value = '1.828.06'
$ws.Range('E12').Value = '  -0.40%  '

# Row 13
$ws.Range('E13').Value = '  -3.80%  '

# Row 14
$ws.Range('D14').Value = '''0.000009905'
$ws.Range('E14').Value = '  -3.53%  '

# Row 15
$ws.Range('D15').Value = '''0.6252'
$ws.Range('E15').Value = '  -7.60%  '

# Row 16
$ws.Range('D16').Value = '2.075.51'
$ws.Range('E16').Value = '  -0.71%  '

# Row 17
$ws.Range('D17').Value = '''78.94'
$ws.Range('E17').Value = '  -3.43%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '29.190.24'
$ws.Range('E18').Value = '  -0.84%  '

# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '''5.826'
$ws.Range('E19').Value = '  -6.68%  '

# Row 20
$ws.Range('D20').Value = '''225.74'
$ws.Range('E20').Value = '  -2.54%  '

# Row 21
$ws.Range('D21').Value = '''1.002'
$ws.Range('E21').Value = '  +0.13%  '

# Row 22
$ws.Range('E22').Value = '  -5.14%  '

# Row 23
$ws.Range('D23').Value = '''6.990'
$ws.Range('E23').Value = '  -4.45%  '

# Row 24
$ws.Range('E24').Value = '  +0.16%  '

# Row 25
$ws.Range('D25').Value = '''155.42'
$ws.Range('E25').Value = '  -1.88%  '

# Row 26
$ws.Range('D26').Value = '''7.997'
$ws.Range('E26').Value = '  -5.72%  '

# Row 27
$ws.Range('E27').Value = '  -4.05%  '

# Row 28
$ws.Range('D28').Value = '''16.54'
$ws.Range('E28').Value = '  -4.76%  '

# Row 29
$ws.Range('D29').Value = '''1.480'
$ws.Range('E29').Value = '  +1.09%  '

# Row 30
$ws.Range('D30').Value = '''0.06181'
$ws.Range('E30').Value = '  -15.59%  '

# Row 31
$ws.Range('D31').Value = '''1.441'
$ws.Range('E31').Value = '  -2.53%  '

# Row 32
$ws.Range('E32').Value = '  -5.68%  '

# Row 33
$ws.Range('D33').Value = '''3.790'
$ws.Range('E33').Value = '  -6.58%  '

# Row 34
$ws.Range('D34').Value = '''1.120'
$ws.Range('E34').Value = '  -1.80%  '

# Row 35
$ws.Range('D35').Value = '''1.744'
$ws.Range('E35').Value = '  -3.95%  '

# Row 36
$ws.Range('D36').Value = '''0.6393'
$ws.Range('E36').Value = '  -8.34%  '

# Row 37
$ws.Range('D37').Value = '''2.537'
$ws.Range('E37').Value = '  -1.11%  '

# Row 38
$ws.Range('D38').Value = '1.216.98'
$ws.Range('E38').Value = '  -1.51%  '

# Row 39
$ws.Range('E39').Value = '  -3.03%  '

# Row 40
$ws.Range('D40').Value = '''0.01741'
$ws.Range('E40').Value = '  -5.30%  '

# Row 41
$ws.Range('D41').Value = '''6.528'
$ws.Range('E41').Value = '  -6.24%  '

# Row 42
$ws.Range('D42').Value = '''0.9010'
$ws.Range('E42').Value = '  -4.73%  '

# Row 43
$ws.Range('D43').Value = '''1.002'
$ws.Range('E43').Value = '  +0.21%  '

# Row 44
$ws.Range('D44').Value = '1.984.36'
$ws.Range('E44').Value = '  -0.56%  '

# Row 45
$ws.Range('D45').Value = '''100.48'
$ws.Range('E45').Value = '  -0.26%  '

# Row 46
$ws.Range('D46').Value = '''62.56'
$ws.Range('E46').Value = '  -4.71%  '

# Row 47
$ws.Range('D47').Value = '''0.00000000117'
$ws.Range('E47').Value = '  -1.83%  '

# Row 48
$ws.Range('D48').Value = '''8.518'
$ws.Range('E48').Value = '  -4.92%  '

# Row 49
$ws.Range('D49').Value = '''1.576'
$ws.Range('E49').Value = '  -8.44%  '

# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.4555'
$ws.Range('E50').Value = '  -0.61%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05513'
$ws.Range('E51').Value = '  -2.58%  '
